# Apply "results with fixed workflow" update to both worksheets (NBR, BAR):
#   - drop the last 4 data rows (rows 17-20), shrinking the table to A1:C16
#   - update the "Cutoff" (col B) and "Reaction_number" (col C) values for the
#     remaining 15 data rows (rows 2-16) to the newly recomputed numbers

$wb = $excel.ActiveWorkbook

# New Cutoff (B) / Reaction_number (C) values per sheet, for rows 2-16 (A stays 0-14).
$nbrValues = @(
  @(5,131),
  @(6,131),
  @(7,130),
  @(8,128),
  @(9,129),
  @(10,127),
  @(11,127),
  @(12,121),
  @(13,120),
  @(14,118),
  @(15,118),
  @(16,118),
  @(17,116),
  @(18,112),
  @(19,112)
)

$barValues = @(
  @(5,607),
  @(6,606),
  @(7,605),
  @(8,608),
  @(9,605),
  @(10,605),
  @(11,605),
  @(12,607),
  @(13,602),
  @(14,600),
  @(15,598),
  @(16,600),
  @(17,599),
  @(18,600),
  @(19,597)
)

$sheetData = @{
  "NBR" = $nbrValues
  "BAR" = $barValues
}

foreach ($ws in $wb.Worksheets) {
  $values = $sheetData[$ws.Name]
  if ($values -eq $null) {
    continue
  }

  # Remove the now-obsolete trailing rows (17-20) first.
  $ws.Range("A17:C20").EntireRow.Delete()

  # Rewrite the Cutoff / Reaction_number columns for the remaining rows.
  for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
  }
}
